$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing logout time for row 13 (D13)
$ws.Cells.Item(13, 4).Value = "2025-07-12 16:07:35"

# Add new row 14 with a fresh login/logout record
$ws.Cells.Item(14, 1).Value = "kumarshashwat890@gmail.com"
$ws.Cells.Item(14, 2).Value = "Shashwat kumar"
$ws.Cells.Item(14, 3).Value = "2025-07-12 16:07:47"
$ws.Cells.Item(14, 4).Value = "2025-07-12 16:09:06"
